$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (well outside the used range) used to produce a pure text
# value via a formula result, then paste-special as a value. This avoids
# Excel's autodetection that would otherwise convert a "YYYY-MM-DD"-looking
# string typed directly into a cell into a date serial number.
$helper = $ws.Cells.Item(1000, 1000)

for ($row = 2; $row -le 31; $row++) {
    $helper.Formula = "=""2013-06-12"""
    $helper.Copy()
    $dst = $ws.Cells.Item($row, 58)
    $dst.PasteSpecial(-4163)
}

$helper.Clear()
$excel.CutCopyMode = $false
